# 16/11/22 Refactoring 2 - Three State Program
# The "Prayer Requests (Image)" slide (the second slide in the deck) is
# being retired in favour of the checkbox-driven Prayer Requests slide
# that already exists later in the deck, so it is simply deleted here.

$p = $ppt.ActivePresentation

# The slide holding the old image-based "Prayer Requests (Image)" title
# is the 2nd slide in the presentation (sldId 362). Removing it shifts
# every following slide up by one position, matching the new sldIdLst
# order (357, 358, 360, 361, 359, 363, 364).
$s = $p.Slides.Item(2)
$s.Delete()
